$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Insert a new column before column B, shifting existing B:I to C:J
$ws.Range("B:B").Insert()

# --- Row 1 headers ---
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "P1"
$ws.Range("C1").Value = "P2"
$ws.Range("D1").Value = "P3"
$ws.Range("E1").Value = "P4"
$ws.Range("F1").Value = "P5"
$ws.Range("G1").Value = "P6"
$ws.Range("H1").Value = "P7"
$ws.Range("I1").Value = "P8"
$ws.Range("J1").Value = "P9"

# --- Row 2 data ---
$ws.Range("A2").Value = "NA"
$ws.Range("B2").Value = "Bacharelado"
$ws.Range("C2").Value = "Ecologia"
$ws.Range("D2").Value = "Não"
$ws.Range("E2").Value = "NA"
$ws.Range("F2").Value = "Não"
$ws.Range("G2").Value = "NA"
$ws.Range("H2").Value = "Não"
$ws.Range("I2").Value = "HTML"
$ws.Range("J2").Value = "PPGE"

# --- Formatting ---
# Whole columns default centered (mirrors target col style)
$ws.Cells.HorizontalAlignment = -4108

# Header row: centered horizontally and vertically
$header = $ws.Range("A1:J1")
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4108

# Data row: centered horizontally (already covered by column default, but explicit for safety)
$ws.Range("A2:J2").HorizontalAlignment = -4108

# Leave the header row selected, matching the saved selection state
$ws.Range("A1:J1").Select() | Out-Null
